$wb = $excel.ActiveWorkbook

# --- Update the Date metadata value on the "Metadata" sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-05-05T11:54:16+00:00"

# --- Update Min/Max/Base Min/Base Max for ActorXDS.XCN9.composant1 (row 7) ---
# on the "Elements" sheet: cardinality changes from 1 to 0.
# A leading apostrophe keeps these as text ("0") instead of Excel
# auto-converting the numeric-looking value to a number, matching the
# original cell type (text stored as a shared string).
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("F7").Value = "'0"
$elements.Range("G7").Value = "'0"
$elements.Range("AG7").Value = "'0"
$elements.Range("AH7").Value = "'0"
